$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 80 (2020-06-05) ---
$ws.Range("A80").Value = 43987
$ws.Range("B80").Formula = "=B79+C80"
$ws.Range("C80").Value = 30
$ws.Range("D80").Formula = "=((B80-B79)/B79)*100"
$ws.Range("E80").Value = 47935
$ws.Range("F80").Formula = "=E80-E79"
$ws.Range("G80").Formula = "=AVERAGE(F78:F80)"
$ws.Range("H80").Formula = "=((E80-E79)/E79)*100"
$ws.Range("I80").Formula = "=B80"
$ws.Range("J80").Formula = "=(I80/E80)*100"
$ws.Range("K80").Formula = "=(C80/F80)*100"
$ws.Range("L80").Value = 168
$ws.Range("M80").Value = 230
$ws.Range("N80").Formula = "=(M80/B80)*100"
$ws.Range("O80").Formula = "=((M80-M79)/M79)*100"
$ws.Range("P80").Value = 16
$ws.Range("Q80").Formula = "=(P80/L80)*100"
$ws.Range("R80").Value = 902
$ws.Range("S80").Formula = "=R80-L80"
$ws.Range("T80").Formula = "=(L80/R80)*100"
$ws.Range("U80").Formula = "=(P80/R80)*100"
$ws.Range("V80").Formula = "=L80-P80"

# --- Row 81 (2020-06-06) ---
$ws.Range("A81").Value = 43988
$ws.Range("B81").Formula = "=B80+C81"
$ws.Range("C81").Value = 32
$ws.Range("D81").Formula = "=((B81-B80)/B80)*100"
$ws.Range("E81").Value = 50621
$ws.Range("F81").Formula = "=E81-E80"
$ws.Range("G81").Formula = "=AVERAGE(F79:F81)"
$ws.Range("H81").Formula = "=((E81-E80)/E80)*100"
$ws.Range("I81").Formula = "=B81"
$ws.Range("J81").Formula = "=(I81/E81)*100"
$ws.Range("K81").Formula = "=(C81/F81)*100"
$ws.Range("L81").Value = 164
$ws.Range("M81").Value = 235
$ws.Range("N81").Formula = "=(M81/B81)*100"
$ws.Range("O81").Formula = "=((M81-M80)/M80)*100"
$ws.Range("P81").Value = 16
$ws.Range("Q81").Formula = "=(P81/L81)*100"
$ws.Range("R81").Value = 851
$ws.Range("S81").Formula = "=R81-L81"
$ws.Range("T81").Formula = "=(L81/R81)*100"
$ws.Range("U81").Formula = "=(P81/R81)*100"
$ws.Range("V81").Formula = "=L81-P81"

# --- Row 82 (2020-06-07) ---
$ws.Range("A82").Value = 43989
$ws.Range("B82").Formula = "=B81+C82"
$ws.Range("C82").Value = 32
$ws.Range("D82").Formula = "=((B82-B81)/B81)*100"
$ws.Range("E82").Value = 51455
$ws.Range("F82").Formula = "=E82-E81"
$ws.Range("G82").Formula = "=AVERAGE(F80:F82)"
$ws.Range("H82").Formula = "=((E82-E81)/E81)*100"
$ws.Range("I82").Formula = "=B82"
$ws.Range("J82").Formula = "=(I82/E82)*100"
$ws.Range("K82").Formula = "=(C82/F82)*100"
$ws.Range("L82").Value = 167
$ws.Range("M82").Value = 236
$ws.Range("N82").Formula = "=(M82/B82)*100"
$ws.Range("O82").Formula = "=((M82-M81)/M81)*100"
$ws.Range("P82").Value = 17
$ws.Range("Q82").Formula = "=(P82/L82)*100"
$ws.Range("R82").Value = 870
$ws.Range("S82").Formula = "=R82-L82"
$ws.Range("T82").Formula = "=(L82/R82)*100"
$ws.Range("U82").Formula = "=(P82/R82)*100"
$ws.Range("V82").Formula = "=L82-P82"

# Column A on the new rows is a date - copy the date format from the last
# existing row (A79) so the new cells render as dates (matches existing style).
$ws.Range("A79").Copy()
$ws.Range("A80:A82").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the author's final selection/viewport from the diff.
$ws.Range("R83").Select()

$wb.Save()
